$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 71429976
$ws.Range("I135").Value = 29413324
$ws.Range("K135").Value = 264719916
$ws.Range("M135").Value = -264717381
$ws.Range("H137").Value = 2357.7234
$ws.Range("I137").Value = 2259.2307
$ws.Range("J137").Value = 2837.875
$ws.Range("K137").Value = 6777.6921
$ws.Range("L137").Value = 8513.625
$ws.Range("M137").Value = -4227.6921
$ws.Range("N137").Value = -13613.625

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 9701.975
$ws.Range("I61").Value = 5318
$ws.Range("J61").Value = 17530.5
$ws.Range("K61").Value = 5318
$ws.Range("L61").Value = 17530.5
$ws.Range("M61").Value = -5106
$ws.Range("N61").Value = -17954.5
$ws.Range("H74").Value = 4992.3145
$ws.Range("I74").Value = 2034.2307
$ws.Range("J74").Value = 13537.889
$ws.Range("K74").Value = 2034.2307
$ws.Range("L74").Value = 13537.889
$ws.Range("M74").Value = -1160.2307
$ws.Range("N74").Value = -15285.889
$ws.Range("H77").Value = 4992.3145
$ws.Range("I77").Value = 2034.2307
$ws.Range("J77").Value = 13537.889
$ws.Range("K77").Value = 10171.1535
$ws.Range("L77").Value = 67689.44499999999
$ws.Range("M77").Value = -5803.1535
$ws.Range("N77").Value = -76425.44499999999
$ws.Range("H122").Value = 1958.375
$ws.Range("I122").Value = 1958.375
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5875.125
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -3425.125
$ws.Range("H136").Value = 9701.975
$ws.Range("I136").Value = 5318
$ws.Range("J136").Value = 17530.5
$ws.Range("K136").Value = 15954
$ws.Range("L136").Value = 52591.5
$ws.Range("M136").Value = -13404
$ws.Range("N136").Value = -57691.5
$ws.Range("H138").Value = 63497.5
$ws.Range("J138").Value = 63497.5
$ws.Range("L138").Value = 63497.5
$ws.Range("N138").Value = -73777.5
$ws.Range("H140").Value = 55921.668
$ws.Range("J140").Value = 55921.668
$ws.Range("L140").Value = 55921.668
$ws.Range("N140").Value = -66281.66800000001
$ws.Range("H141").Value = 53600
$ws.Range("J141").Value = 53600
$ws.Range("L141").Value = 53600
$ws.Range("N141").Value = -63960

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 32139.234
$ws.Range("I134").Value = 2545.1853
$ws.Range("J134").Value = 146287.72
$ws.Range("K134").Value = 7635.5559
$ws.Range("L134").Value = 438863.16
$ws.Range("M134").Value = -5100.5559
$ws.Range("N134").Value = -443933.16

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6889.7954
$ws.Range("I31").Value = 6677.528
$ws.Range("J31").Value = 7845
$ws.Range("K31").Value = 6677.528
$ws.Range("L31").Value = 7845
$ws.Range("M31").Value = -6382.528
$ws.Range("N31").Value = -8435
$ws.Range("H34").Value = 6889.7954
$ws.Range("I34").Value = 6677.528
$ws.Range("J34").Value = 7845
$ws.Range("K34").Value = 6677.528
$ws.Range("L34").Value = 7845
$ws.Range("M34").Value = -6475.528
$ws.Range("N34").Value = -8249
$ws.Range("H58").Value = 2846299.5
$ws.Range("I58").Value = 4330986.5
$ws.Range("J58").Value = 11897.091
$ws.Range("K58").Value = 4330986.5
$ws.Range("L58").Value = 11897.091
$ws.Range("M58").Value = -4330783.5
$ws.Range("N58").Value = -12303.091
$ws.Range("H122").Value = 12330
$ws.Range("I122").Value = 11996
$ws.Range("J122").Value = 14000
$ws.Range("K122").Value = 35988
$ws.Range("L122").Value = 42000
$ws.Range("M122").Value = -33538
$ws.Range("N122").Value = -46900
$ws.Range("H132").Value = 5036.7617
$ws.Range("I132").Value = 5581.7334
$ws.Range("J132").Value = 3674.3333
$ws.Range("K132").Value = 16745.2002
$ws.Range("L132").Value = 11022.9999
$ws.Range("M132").Value = -14215.2002
$ws.Range("N132").Value = -16082.9999
$ws.Range("H134").Value = 2322.348
$ws.Range("I134").Value = 1818.0605
$ws.Range("J134").Value = 3602.4614
$ws.Range("K134").Value = 5454.181500000001
$ws.Range("L134").Value = 10807.3842
$ws.Range("M134").Value = -2919.181500000001
$ws.Range("N134").Value = -15877.3842
$ws.Range("H136").Value = 2846299.5
$ws.Range("I136").Value = 4330986.5
$ws.Range("J136").Value = 11897.091
$ws.Range("K136").Value = 12992959.5
$ws.Range("L136").Value = 35691.273
$ws.Range("M136").Value = -12990409.5
$ws.Range("N136").Value = -40791.273

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 27003
$ws.Range("I122").Value = 27003
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 81009
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -78559
$ws.Range("H126").Value = 2683.1304
$ws.Range("I126").Value = 1885.5385
$ws.Range("J126").Value = 3720
$ws.Range("K126").Value = 5656.6155
$ws.Range("L126").Value = 11160
$ws.Range("M126").Value = -3186.6155
$ws.Range("N126").Value = -16100

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5347.1177
$ws.Range("I132").Value = 5275.5
$ws.Range("J132").Value = 5410.778
$ws.Range("K132").Value = 15826.5
$ws.Range("L132").Value = 16232.334
$ws.Range("M132").Value = -13296.5
$ws.Range("N132").Value = -21292.334
$ws.Range("H136").Value = 6080.5483
$ws.Range("I136").Value = 3506.9285
$ws.Range("K136").Value = 10520.7855
$ws.Range("M136").Value = -7970.7855

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 18367.5
$ws.Range("I122").Value = 2100
$ws.Range("J122").Value = 26501.25
$ws.Range("K122").Value = 6300
$ws.Range("L122").Value = 79503.75
$ws.Range("M122").Value = -3850
$ws.Range("N122").Value = -84403.75
$ws.Range("H132").Value = 3267.4167
$ws.Range("I132").Value = 2420.5
$ws.Range("J132").Value = 7502
$ws.Range("K132").Value = 7261.5
$ws.Range("L132").Value = 22506
$ws.Range("M132").Value = -4731.5
$ws.Range("N132").Value = -27566
$ws.Range("H136").Value = 6419.615
$ws.Range("I136").Value = 3186.4736
$ws.Range("J136").Value = 9491.1
$ws.Range("K136").Value = 9559.4208
$ws.Range("L136").Value = 28473.3
$ws.Range("M136").Value = -7009.4208
$ws.Range("N136").Value = -33573.3

Write-Output "All edits applied"